$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.143.45"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "1.679.91"
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'215.25"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +1.93%  "
$ws.Range("D9").Value = "'21.35"
$ws.Range("E9").Value = "  +5.11%  "
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").Value = "'0.0886"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.917.44"
$ws.Range("D13").Value = "1.672.14"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("E14").Value = "  +1.41%  "
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "'66.25"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "27.152.00"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "'239.27"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "'8.06"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D20").Value = "0.0₃0742"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "'9.44"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").Value = "'147.00"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "'7.24"
$ws.Range("E26").Value = "  +0.84%  "
$ws.Range("E27").Value = "  +2.14%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "1.562.85"
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("E34").Value = "  +2.83%  "
$ws.Range("D35").Value = "'1.68"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "'0.601"
$ws.Range("E36").Value = "  +2.78%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  +4.18%  "
$ws.Range("D39").Value = "'0.0174"
$ws.Range("E39").Value = "  +2.58%  "
$ws.Range("E40").Value = "  +1.24%  "
$ws.Range("D41").Value = "'69.15"
$ws.Range("E41").Value = "  +2.70%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "'5.58"
$ws.Range("E43").Value = "  -4.67%  "
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").Value = "1.826.16"
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("D46").Value = "'0.785"
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").Value = "'90.62"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").Value = "'1.59"
$ws.Range("E48").Value = "  +3.71%  "
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("E50").Value = "  +2.29%  "
$ws.Range("E51").Value = "  +6.19%  "
